# Update the "veicolo" column (N) values on the "Schedulazione" sheet,
# swapping the " (esterno)" / " (non in estrazione)" suffix for a set of rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Schedulazione")

$ws.Range("N5").Value  = "39666 (esterno)"
$ws.Range("N6").Value  = "39742 (non in estrazione)"
$ws.Range("N11").Value = "39666 (non in estrazione)"
$ws.Range("N20").Value = "39762 (non in estrazione)"
$ws.Range("N21").Value = "39723 (esterno)"
$ws.Range("N26").Value = "39750 (esterno)"
$ws.Range("N27").Value = "39764 (esterno)"
